$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("stok") held strings mirroring the "satuan" column by mistake.
# Replace them with actual numeric stock values.
$ws.Range("D2").Value = 10
$ws.Range("D3").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("D5").Value = 0
$ws.Range("D6").Value = 0
$ws.Range("D7").Value = 0
